$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.059880695826387
$ws.Cells.Item(2, 4).Value = 1.066195651312385
$ws.Cells.Item(2, 5).Value = 1.06559725020907
$ws.Cells.Item(2, 6).Value = 1.076294384897704
$ws.Cells.Item(2, 9).Value = 1.02359499962809
$ws.Cells.Item(2, 10).Value = 1.064864431662581
$ws.Cells.Item(2, 11).Value = 1.068906852217692
$ws.Cells.Item(2, 12).Value = 1.06831006178296
$ws.Cells.Item(2, 13).Value = 1.078978695818275
$ws.Cells.Item(2, 14).Value = 1.024944988265652

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.061845146055405
$ws.Cells.Item(3, 4).Value = 1.068008479248407
$ws.Cells.Item(3, 5).Value = 1.067371612228636
$ws.Cells.Item(3, 6).Value = 1.07820486085402
$ws.Cells.Item(3, 9).Value = 1.023504579208684
$ws.Cells.Item(3, 10).Value = 1.066477460189505
$ws.Cells.Item(3, 11).Value = 1.070532738703228
$ws.Cells.Item(3, 12).Value = 1.069897459242697
$ws.Cells.Item(3, 13).Value = 1.080703977836677
$ws.Cells.Item(3, 14).Value = 1.025517505399298

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.063113377602309
$ws.Cells.Item(4, 4).Value = 1.069178901813843
$ws.Cells.Item(4, 5).Value = 1.068517146255217
$ws.Cells.Item(4, 6).Value = 1.079438591075603
$ws.Cells.Item(4, 9).Value = 1.023443335729026
$ws.Cells.Item(4, 10).Value = 1.067518081017153
$ws.Cells.Item(4, 11).Value = 1.071581758037445
$ws.Cells.Item(4, 12).Value = 1.070921567324904
$ws.Cells.Item(4, 13).Value = 1.081817449544657
$ws.Cells.Item(4, 14).Value = 1.025886005745347

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.063645868060701
$ws.Cells.Item(5, 4).Value = 1.06967034282929
$ws.Cells.Item(5, 5).Value = 1.068998124028754
$ws.Cells.Item(5, 6).Value = 1.079956676848268
$ws.Cells.Item(5, 9).Value = 1.023416932628374
$ws.Cells.Item(5, 10).Value = 1.067954828209848
$ws.Cells.Item(5, 11).Value = 1.072022054706459
$ws.Cells.Item(5, 12).Value = 1.071351389393252
$ws.Cells.Item(5, 13).Value = 1.08228487548647
$ws.Cells.Item(5, 14).Value = 1.026040460367392

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.063735236599872
$ws.Cells.Item(6, 4).Value = 1.0697528229853
$ws.Cells.Item(6, 5).Value = 1.069078847341389
$ws.Cells.Item(6, 6).Value = 1.080043632571595
$ws.Cells.Item(6, 9).Value = 1.023412460900047
$ws.Cells.Item(6, 10).Value = 1.068028117659704
$ws.Cells.Item(6, 11).Value = 1.072095941226434
$ws.Cells.Item(6, 12).Value = 1.071423517057895
$ws.Cells.Item(6, 13).Value = 1.08236331911628
$ws.Cells.Item(6, 14).Value = 1.026066366989902

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.06312049538743
$ws.Cells.Item(7, 4).Value = 1.069185470824823
$ws.Cells.Item(7, 5).Value = 1.068523575456853
$ws.Cells.Item(7, 6).Value = 1.079445515995302
$ws.Cells.Item(7, 9).Value = 1.02344298551034
$ws.Cells.Item(7, 10).Value = 1.06752391969997
$ws.Cells.Item(7, 11).Value = 1.071587644074963
$ws.Cells.Item(7, 12).Value = 1.070927313408352
$ws.Cells.Item(7, 13).Value = 1.081823697952229
$ws.Cells.Item(7, 14).Value = 1.025888071387038

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.060545203331269
$ws.Cells.Item(8, 4).Value = 1.06680885353093
$ws.Cells.Item(8, 5).Value = 1.066197452109527
$ws.Cells.Item(8, 6).Value = 1.076940562770609
$ws.Cells.Item(8, 9).Value = 1.023565007130014
$ws.Cells.Item(8, 10).Value = 1.065410217344836
$ws.Cells.Item(8, 11).Value = 1.069456967007026
$ws.Cells.Item(8, 12).Value = 1.068847170079768
$ws.Cells.Item(8, 13).Value = 1.079562373911781
$ws.Cells.Item(8, 14).Value = 1.02513888153306

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.05598410572061
$ws.Cells.Item(9, 4).Value = 1.062600264868534
$ws.Cells.Item(9, 5).Value = 1.062077872050329
$ws.Cells.Item(9, 6).Value = 1.072506734448362
$ws.Cells.Item(9, 9).Value = 1.023759139049229
$ws.Cells.Item(9, 10).Value = 1.06166100050086
$ws.Cells.Item(9, 11).Value = 1.065678442843346
$ws.Cells.Item(9, 12).Value = 1.06515766306946
$ws.Cells.Item(9, 13).Value = 1.075554647605434
$ws.Cells.Item(9, 14).Value = 1.023803481430824

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.05292657621234
$ws.Cells.Item(10, 4).Value = 1.059779528477301
$ws.Cells.Item(10, 5).Value = 1.059316524790412
$ws.Cells.Item(10, 6).Value = 1.069536395660843
$ws.Cells.Item(10, 9).Value = 1.023874603872166
$ws.Cells.Item(10, 10).Value = 1.059143948768712
$ws.Cells.Item(10, 11).Value = 1.063142262967229
$ws.Cells.Item(10, 12).Value = 1.062680833646136
$ws.Cells.Item(10, 13).Value = 1.072866301648161
$ws.Cells.Item(10, 14).Value = 1.022902625331242

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.051598365262403
$ws.Cells.Item(11, 4).Value = 1.058554311315541
$ws.Cells.Item(11, 5).Value = 1.05811704204199
$ws.Cells.Item(11, 6).Value = 1.068246519868483
$ws.Cells.Item(11, 9).Value = 1.023921309042974
$ws.Cells.Item(11, 10).Value = 1.058049638298319
$ws.Cells.Item(11, 11).Value = 1.062039767897639
$ws.Cells.Item(11, 12).Value = 1.061604045054462
$ws.Cells.Item(11, 13).Value = 1.07169805724299
$ws.Cells.Item(11, 14).Value = 1.022509954756825

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.051104341028799
$ws.Cells.Item(12, 4).Value = 1.05809861587819
$ws.Cells.Item(12, 5).Value = 1.057670908844888
$ws.Cells.Item(12, 6).Value = 1.067766824418723
$ws.Cells.Item(12, 9).Value = 1.023938164268129
$ws.Cells.Item(12, 10).Value = 1.057642480039465
$ws.Cells.Item(12, 11).Value = 1.06162958431155
$ws.Cells.Item(12, 12).Value = 1.061203411594794
$ws.Cells.Item(12, 13).Value = 1.071263471498063
$ws.Cells.Item(12, 14).Value = 1.022363702895156

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.051210341507673
$ws.Cells.Item(13, 4).Value = 1.058196391377689
$ws.Cells.Item(13, 5).Value = 1.057766633051658
$ws.Cells.Item(13, 6).Value = 1.067869747245153
$ws.Cells.Item(13, 9).Value = 1.023934571058996
$ws.Cells.Item(13, 10).Value = 1.0577298480822
$ws.Cells.Item(13, 11).Value = 1.061717600626473
$ws.Cells.Item(13, 12).Value = 1.061289379294981
$ws.Cells.Item(13, 13).Value = 1.071356721265141
$ws.Cells.Item(13, 14).Value = 1.022395092470031

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.051557542825475
$ws.Cells.Item(14, 4).Value = 1.058516655655621
$ws.Cells.Item(14, 5).Value = 1.05808017672863
$ws.Cells.Item(14, 6).Value = 1.068206880017272
$ws.Cells.Item(14, 9).Value = 1.023922712353276
$ws.Cells.Item(14, 10).Value = 1.058015996511619
$ws.Cells.Item(14, 11).Value = 1.062005875734166
$ws.Cells.Item(14, 12).Value = 1.061570942271987
$ws.Cells.Item(14, 13).Value = 1.071662147585233
$ws.Cells.Item(14, 14).Value = 1.022497873661107

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.051771375823082
$ws.Cells.Item(15, 4).Value = 1.05871390152969
$ws.Cells.Item(15, 5).Value = 1.058273282276895
$ws.Cells.Item(15, 6).Value = 1.068414521361465
$ws.Cells.Item(15, 9).Value = 1.023915340507297
$ws.Cells.Item(15, 10).Value = 1.058192210945733
$ws.Cells.Item(15, 11).Value = 1.062183402447911
$ws.Cells.Item(15, 12).Value = 1.061744333646821
$ws.Cells.Item(15, 13).Value = 1.071850244359645
$ws.Cells.Item(15, 14).Value = 1.022561147795983

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.053014632933413
$ws.Cells.Item(16, 4).Value = 1.059860759811118
$ws.Cells.Item(16, 5).Value = 1.059396048653099
$ws.Cells.Item(16, 6).Value = 1.069621920624346
$ws.Cells.Item(16, 9).Value = 1.023871434968932
$ws.Cells.Item(16, 10).Value = 1.059216479925519
$ws.Cells.Item(16, 11).Value = 1.063215339380079
$ws.Cells.Item(16, 12).Value = 1.062752204201316
$ws.Cells.Item(16, 13).Value = 1.072943744492237
$ws.Cells.Item(16, 14).Value = 1.022928630352987

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.053793332094963
$ws.Cells.Item(17, 4).Value = 1.060579116072283
$ws.Cells.Item(17, 5).Value = 1.06009929804501
$ws.Cells.Item(17, 6).Value = 1.070378284137798
$ws.Cells.Item(17, 9).Value = 1.023843013862359
$ws.Cells.Item(17, 10).Value = 1.059857781946441
$ws.Cells.Item(17, 11).Value = 1.063861477596405
$ws.Cells.Item(17, 12).Value = 1.063383248526231
$ws.Cells.Item(17, 13).Value = 1.073628536416646
$ws.Cells.Item(17, 14).Value = 1.023158443197648

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.054247122801946
$ws.Cells.Item(18, 4).Value = 1.060997754037345
$ws.Cells.Item(18, 5).Value = 1.060509126179034
$ws.Cells.Item(18, 6).Value = 1.070819102753552
$ws.Cells.Item(18, 9).Value = 1.023826118601784
$ws.Cells.Item(18, 10).Value = 1.060231418383327
$ws.Cells.Item(18, 11).Value = 1.064237944356233
$ws.Cells.Item(18, 12).Value = 1.063750911949206
$ws.Cells.Item(18, 13).Value = 1.074027562759262
$ws.Cells.Item(18, 14).Value = 1.023292239329578

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.054401784671152
$ws.Cells.Item(19, 4).Value = 1.061140437142805
$ws.Cells.Item(19, 5).Value = 1.060648805640096
$ws.Cells.Item(19, 6).Value = 1.070969350895389
$ws.Cells.Item(19, 9).Value = 1.023820303822167
$ws.Cells.Item(19, 10).Value = 1.060358747405628
$ws.Cells.Item(19, 11).Value = 1.064366240046567
$ws.Cells.Item(19, 12).Value = 1.06387620602876
$ws.Cells.Item(19, 13).Value = 1.074163553006096
$ws.Cells.Item(19, 14).Value = 1.023337818185289

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.053709827767501
$ws.Cells.Item(20, 4).Value = 1.06050208136764
$ws.Cells.Item(20, 5).Value = 1.060023883979475
$ws.Cells.Item(20, 6).Value = 1.070297170415311
$ws.Cells.Item(20, 9).Value = 1.023846096030143
$ws.Cells.Item(20, 10).Value = 1.059789020328171
$ws.Cells.Item(20, 11).Value = 1.063792196121958
$ws.Cells.Item(20, 12).Value = 1.06331558641015
$ws.Cells.Item(20, 13).Value = 1.073555106372267
$ws.Cells.Item(20, 14).Value = 1.023133812342748

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.051455319316625
$ws.Cells.Item(21, 4).Value = 1.058422362412335
$ws.Cells.Item(21, 5).Value = 1.057987862438356
$ws.Cells.Item(21, 6).Value = 1.068107618959234
$ws.Cells.Item(21, 9).Value = 1.02392621804881
$ws.Cells.Item(21, 10).Value = 1.057931751936927
$ws.Cells.Item(21, 11).Value = 1.061921004496668
$ws.Cells.Item(21, 12).Value = 1.061488047541635
$ws.Cells.Item(21, 13).Value = 1.071572225196355
$ws.Cells.Item(21, 14).Value = 1.022467618163366

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.050033945810125
$ws.Cells.Item(22, 4).Value = 1.057111305962144
$ws.Cells.Item(22, 5).Value = 1.056704299849443
$ws.Cells.Item(22, 6).Value = 1.066727604795404
$ws.Cells.Item(22, 9).Value = 1.023973741382756
$ws.Cells.Item(22, 10).Value = 1.056760054029482
$ws.Cells.Item(22, 11).Value = 1.060740638148382
$ws.Cells.Item(22, 12).Value = 1.060335136896347
$ws.Cells.Item(22, 13).Value = 1.070321750871319
$ws.Cells.Item(22, 14).Value = 1.022046457816724

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.050787818305922
$ws.Cells.Item(23, 4).Value = 1.057806656451157
$ws.Cells.Item(23, 5).Value = 1.057385073282811
$ws.Cells.Item(23, 6).Value = 1.067459502080612
$ws.Cells.Item(23, 9).Value = 1.023948818288664
$ws.Cells.Item(23, 10).Value = 1.057381575323464
$ws.Cells.Item(23, 11).Value = 1.061366746604083
$ws.Cells.Item(23, 12).Value = 1.060946689499404
$ws.Cells.Item(23, 13).Value = 1.070985014230036
$ws.Cells.Item(23, 14).Value = 1.022269942973165

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.053747561049086
$ws.Cells.Item(24, 4).Value = 1.060536891168557
$ws.Cells.Item(24, 5).Value = 1.060057961477949
$ws.Cells.Item(24, 6).Value = 1.07033382331275
$ws.Cells.Item(24, 9).Value = 1.023844704313072
$ws.Cells.Item(24, 10).Value = 1.059820092053894
$ws.Cells.Item(24, 11).Value = 1.063823502719114
$ws.Cells.Item(24, 12).Value = 1.063346161286895
$ws.Cells.Item(24, 13).Value = 1.07358828748778
$ws.Cells.Item(24, 14).Value = 1.023144942737654

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.057166125959337
$ws.Cells.Item(25, 4).Value = 1.063690850005483
$ws.Cells.Item(25, 5).Value = 1.063145439626395
$ws.Cells.Item(25, 6).Value = 1.07365544743386
$ws.Cells.Item(25, 9).Value = 1.02371141742819
$ws.Cells.Item(25, 10).Value = 1.062633284044241
$ws.Cells.Item(25, 11).Value = 1.066658231033259
$ws.Cells.Item(25, 12).Value = 1.066114439089632
$ws.Cells.Item(25, 13).Value = 1.076593574600589
$ws.Cells.Item(25, 14).Value = 1.024150555004768
